$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 5: "ggplot2 " + "effektiv bei Abschlussarbeiten einsetzen" -> one run ---
$para5 = $tr.Paragraphs(5, 1)
$prefix5 = $para5.Characters(1, 8)
$prefix5.Delete() | Out-Null

$para5b = $tr.Paragraphs(5, 1)
$rest5 = $para5b.Characters(1, $para5b.Text.Length - 1)
$rest5.InsertBefore("ggplot2 ") | Out-Null

# --- Paragraph 6: "Aufzeigen welche Möglichkeiten " + "ggplot" + " mit sich bringt" -> one run ---
$para6 = $tr.Paragraphs(6, 1)
$tail6 = $para6.Characters(32, $para6.Text.Length - 31)
$tail6.Delete() | Out-Null

$para6b = $tr.Paragraphs(6, 1)
$head6 = $para6b.Characters(1, $para6b.Text.Length - 1)
$head6.InsertAfter("ggplot2 mit sich bringt") | Out-Null
